$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is referenced from the Overview sheet (B/C columns)
#    and from the per-language sheets (C column), so update every cell that
#    carries it.
# ---------------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now has a real timestamp instead of
#    the zero-date placeholder.
# ---------------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-20 05:20:31"
$wsZhCn.Range("H3").Value = "2016-03-20 05:20:31"

$wsDeDe.Range("H2").Value = "2016-03-20 05:20:46"
$wsDeDe.Range("H3").Value = "2016-03-20 05:20:46"

# ---------------------------------------------------------------------------
# 3. Populate "Latest Target File" (F) and "Latest Handback File" (G) columns
#    with hyperlinked file names, for both the zh-cn and de-de rows.
# ---------------------------------------------------------------------------

# zh-cn sheet, row 2 (file 1b77ce31...)
$wsZhCn.Range("F2").Value = "1b77ce31-bea1-4e70-9ad4-e555ac224593.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/79e4bcaccf44aa5664403536392064946d42dfec/e2e/1b77ce31-bea1-4e70-9ad4-e555ac224593.md", "", "", "1b77ce31-bea1-4e70-9ad4-e555ac224593.md")

$wsZhCn.Range("G2").Value = "1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/55f306e4585b0fb12aa7225a2364875bc4752e4d/ol-handback/OpenLocalizationTest/oltest/yuwzho/1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.zh-cn.xlf", "", "", "1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.zh-cn.xlf")

# zh-cn sheet, row 3 (file cf10807a...)
$wsZhCn.Range("F3").Value = "cf10807a-6211-4b89-a29d-faca009f048a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/79e4bcaccf44aa5664403536392064946d42dfec/e2e/cf10807a-6211-4b89-a29d-faca009f048a.md", "", "", "cf10807a-6211-4b89-a29d-faca009f048a.md")

$wsZhCn.Range("G3").Value = "cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/55f306e4585b0fb12aa7225a2364875bc4752e4d/ol-handback/OpenLocalizationTest/oltest/yuwzho/cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.zh-cn.xlf", "", "", "cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.zh-cn.xlf")

# de-de sheet, row 2
$wsDeDe.Range("F2").Value = "1b77ce31-bea1-4e70-9ad4-e555ac224593.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/79e4bcaccf44aa5664403536392064946d42dfec/e2e/1b77ce31-bea1-4e70-9ad4-e555ac224593.md", "", "", "1b77ce31-bea1-4e70-9ad4-e555ac224593.md")

$wsDeDe.Range("G2").Value = "1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6364e7ee12e5c299ad2e82cba18fa38085e3aaeb/ol-handback/OpenLocalizationTest/oltest/yuwzho/1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.de-de.xlf", "", "", "1b77ce31-bea1-4e70-9ad4-e555ac224593.1c32694baddcfd3141f899e23bfe4a9b6fe361e2.de-de.xlf")

# de-de sheet, row 3
$wsDeDe.Range("F3").Value = "cf10807a-6211-4b89-a29d-faca009f048a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/79e4bcaccf44aa5664403536392064946d42dfec/e2e/cf10807a-6211-4b89-a29d-faca009f048a.md", "", "", "cf10807a-6211-4b89-a29d-faca009f048a.md")

$wsDeDe.Range("G3").Value = "cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6364e7ee12e5c299ad2e82cba18fa38085e3aaeb/ol-handback/OpenLocalizationTest/oltest/yuwzho/cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.de-de.xlf", "", "", "cf10807a-6211-4b89-a29d-faca009f048a.22efa02f70463d67ed03a0be769ada36ab7c2afe.de-de.xlf")
